$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Fill in the season record values for each data row
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 90   # AD
    $ws.Cells.Item($r, 31).Value = 72   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
